$d = $word.ActiveDocument
$d.Tables.Item(1).Cell(1, 1).Range.Text = "69÷6=11, 3"
$d.Tables.Item(1).Cell(1, 2).Range.Text = "24÷2=12, 0"
$d.Tables.Item(1).Cell(1, 3).Range.Text = "67÷3=22, 1"
$d.Tables.Item(1).Cell(1, 4).Range.Text = "78÷4=19, 2"
$d.Tables.Item(1).Cell(1, 5).Range.Text = "50÷2=25, 0"
$d.Tables.Item(1).Cell(5, 1).Range.Text = "87÷7=12, 3"
$d.Tables.Item(1).Cell(5, 2).Range.Text = "80÷3=26, 2"
$d.Tables.Item(1).Cell(5, 3).Range.Text = "18÷3=6, 0"
$d.Tables.Item(1).Cell(5, 4).Range.Text = "71÷7=10, 1"
$d.Tables.Item(1).Cell(5, 5).Range.Text = "13÷6=2, 1"
$d.Tables.Item(1).Cell(9, 1).Range.Text = "81÷4=20, 1"
$d.Tables.Item(1).Cell(9, 2).Range.Text = "88÷5=17, 3"
$d.Tables.Item(1).Cell(9, 3).Range.Text = "58÷7=8, 2"
$d.Tables.Item(1).Cell(9, 4).Range.Text = "48÷6=8, 0"
$d.Tables.Item(1).Cell(9, 5).Range.Text = "67÷8=8, 3"
$d.Tables.Item(1).Cell(13, 1).Range.Text = "65÷9=7, 2"
$d.Tables.Item(1).Cell(13, 2).Range.Text = "42÷6=7, 0"
$d.Tables.Item(1).Cell(13, 3).Range.Text = "55÷2=27, 1"
$d.Tables.Item(1).Cell(13, 4).Range.Text = "47÷3=15, 2"
$d.Tables.Item(1).Cell(13, 5).Range.Text = "58÷3=19, 1"
$d.Tables.Item(1).Cell(17, 1).Range.Text = "30÷9=3, 3"
$d.Tables.Item(1).Cell(17, 2).Range.Text = "87÷8=10, 7"
$d.Tables.Item(1).Cell(17, 3).Range.Text = "32÷7=4, 4"
$d.Tables.Item(1).Cell(17, 4).Range.Text = "84÷2=42, 0"
$d.Tables.Item(1).Cell(17, 5).Range.Text = "71÷9=7, 8"
